# aggiunta tessera sanitaria + mirror fix
# Appends the new daily COVID-19 data rows (days 97-110, sheet rows 99-112)
# to the existing data table on Sheet1 (columns: Giorno, Infetti, Decessi, Ricoverati).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Giorno, Infetti, Decessi, Ricoverati
$data = @(
    @(97, 201505, 27359, 68941),
    @(98, 203591, 27682, 71252),
    @(99, 205463, 27967, 75945),
    @(100, 207428, 28236, 78249),
    @(101, 209328, 28710, 79914),
    @(102, 210717, 28884, 81654),
    @(103, 211938, 29079, 82879),
    @(104, 213013, 29315, 85231),
    @(105, 214457, 29684, 93245),
    @(106, 215858, 29958, 96276),
    @(107, 217185, 30201, 99023),
    @(108, 218268, 30395, 103031),
    @(109, 219070, 30560, 105186),
    @(110, 219814, 30739, 106587)
)

$startRow = 99
$r = $startRow
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
    $r = $r + 1
}
